$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4046
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 112
$ws.Range("E2").Value = 268
$ws.Range("F2").Value = 17
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1775
$ws.Range("I2").Value = 22264
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 135
$ws.Range("L2").Value = 14456

$ws.Range("B3").Value = 22
$ws.Range("D3").Value = 238
$ws.Range("E3").Value = 238
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 7229
$ws.Range("J3").Value = 5

$ws.Range("D4").Value = 31
$ws.Range("E4").Value = 32
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 6129

$ws.Range("B5").Value = 804
$ws.Range("D5").Value = 119
$ws.Range("E5").Value = 119
$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 8559
$ws.Range("J5").Value = 2

$ws.Range("D6").Value = 49
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 5714
$ws.Range("J6").Value = 0

$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 74
$ws.Range("E8").Value = 96
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 638
$ws.Range("I8").Value = 2297

$ws.Range("B9").Value = 776
$ws.Range("D9").Value = 15
$ws.Range("E9").Value = 15
$ws.Range("I9").Value = 6000

$ws.Range("B10").Value = 774
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 14
$ws.Range("F10").Value = 1
$ws.Range("H10").Value = 3333
$ws.Range("I10").Value = 18333
$ws.Range("J10").Value = 0

$ws.Range("B11").Value = 1002
$ws.Range("D11").Value = 72
$ws.Range("E11").Value = 73
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 20
$ws.Range("I11").Value = 57777
$ws.Range("J11").Value = 4

$ws.Range("B12").Value = 840
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 30000

$ws.Range("B13").Value = 253
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 61
$ws.Range("E13").Value = 101
$ws.Range("F13").Value = 29
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 3187
$ws.Range("I13").Value = 5519
$ws.Range("J13").Value = 3

$ws.Range("B14").Value = 805
$ws.Range("D14").Value = 91
$ws.Range("E14").Value = 92
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 109
$ws.Range("I14").Value = 9000
$ws.Range("J14").Value = 1

$ws.Range("D15").Value = 26
$ws.Range("E15").Value = 27
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 370
$ws.Range("I15").Value = 4231

$ws.Range("B16").Value = 213
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 131
$ws.Range("E16").Value = 183
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 2180
$ws.Range("I16").Value = 12336
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 38
$ws.Range("L16").Value = 2901

$ws.Range("B17").Value = 64
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 3333
$ws.Range("L17").Value = 1429

$ws.Range("B18").Value = 8
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 143
$ws.Range("E18").Value = 178
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 618
$ws.Range("I18").Value = 5734
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 23
$ws.Range("L18").Value = 1292

$ws.Range("B19").Value = 780
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 59
$ws.Range("E19").Value = 63
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 6330

$ws.Range("B20").Value = 1401
$ws.Range("D20").Value = 60
$ws.Range("E20").Value = 78
$ws.Range("F20").Value = 18
$ws.Range("H20").Value = 8800
$ws.Range("I20").Value = 14576
$ws.Range("J20").Value = 0

$ws.Range("B21").Value = 484
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 17
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 25
$ws.Range("I21").Value = 11111
$ws.Range("J21").Value = 0

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0

$ws.Range("B23").Value = 814
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 4
$ws.Range("I23").Value = 2500

$ws.Range("B24").Value = 996
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 17
$ws.Range("H24").Value = 20
$ws.Range("I24").Value = 10357
$ws.Range("J24").Value = 0

$ws.Range("B25").Value = 76
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 3913
$ws.Range("J25").Value = 3
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 8

$ws.Range("B26").Value = 135
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0

$ws.Range("B27").Value = 188
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 33
$ws.Range("E27").Value = 37
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 833
$ws.Range("I27").Value = 2424
$ws.Range("J27").Value = 0
